$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("attributes")

# Insert a new column before column S (shifts mappedBy.. right by one) and
# give the new column the header "maxLength" (exposing maxLength in EMX).
$ws.Columns("S").Insert() | Out-Null
$ws.Range("S1").Value = "maxLength"

# Make "attributes" the active sheet, with the freshly inserted column
# selected (matches the author's post-edit UI state).
$ws.Activate() | Out-Null
$ws.Columns("S").Select() | Out-Null
